# Update odds for the row-3 match (ZED - Al Ahly)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 5.6
$ws.Range("H3").Value = 3.65
$ws.Range("J3").Value = 5.6
$ws.Range("K3").Value = 2.18
$ws.Range("L3").Value = 2.12
$ws.Range("M3").Value = 1.07
$ws.Range("N3").Value = 7
$ws.Range("O3").Value = 1.32
$ws.Range("P3").Value = 3.1
$ws.Range("Q3").Value = 1.98
$ws.Range("R3").Value = 1.78
$ws.Range("S3").Value = 1.39
$ws.Range("T3").Value = 2.75
$ws.Range("U3").Value = 1.98
$ws.Range("V3").Value = 1.75
$ws.Range("X3").Value = 35
$ws.Range("Y3").Value = 18
$ws.Range("Z3").Value = 120
$ws.Range("AC3").Value = 7
$ws.Range("AD3").Value = 7.1
$ws.Range("AG3").Value = 800
$ws.Range("AH3").Value = 6
$ws.Range("AI3").Value = 6.8
$ws.Range("AJ3").Value = 8.25
$ws.Range("AL3").Value = 13.5
$ws.Range("AM3").Value = 30
$ws.Range("AN3").Value = 7.2
$ws.Range("AP3").Value = 35
$ws.Range("AT3").Value = 2.75
$ws.Range("AU3").Value = 7.7
$ws.Range("AV3").Value = 70
$ws.Range("AX3").Value = 7.6
$ws.Range("AY3").Value = 17.5
$ws.Range("AZ3").Value = 25
$ws.Range("BA3").Value = 60
$ws.Range("BB3").Value = 250

# Update odds for the row-5 match (Ruch Chorzow - Chrobry Glogow)
$ws.Range("G5").Value = 1.5
$ws.Range("H5").Value = 4
$ws.Range("I5").Value = 6
$ws.Range("J5").Value = 2.05
$ws.Range("M5").Value = 1.04
$ws.Range("N5").Value = 13
$ws.Range("O5").Value = 1.25
$ws.Range("P5").Value = 3.75
$ws.Range("W5").Value = 7
$ws.Range("AA5").Value = 12
$ws.Range("AC5").Value = 12
$ws.Range("AD5").Value = 8
$ws.Range("AG5").Value = 301
$ws.Range("AI5").Value = 34
$ws.Range("AQ5").Value = 21
$ws.Range("AX5").Value = 34

# Remove the row-6 match (Al Orubah - Al Ittihad) entirely
$ws.Rows(6).Delete()
